# Generate Report for handback
# The two localized source files (b346fae3... and bb845ad7...) have both
# now been handed back and are in sync with en-US. This swaps their
# reporting order (b346fae3 now first, bb845ad7 now second) on every
# sheet, updates the status text, fills in the "Latest Handback File"
# columns with the returned .xlf files and refreshes the "Latest Handback
# DateTime" timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-28 09:51:12"
$ws2.Range("E2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$ws2.Range("F2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-01-28 09:52:02"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-01-28 09:51:12"
$ws2.Range("E3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$ws2.Range("F3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-01-28 09:52:02"
$ws2.Range("H3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27140b88353bd0c13c367205d9e1d43db25ace5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0b837e76f278922680c4f162175a90d15a1c438c/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fd7690f6242f2ab8aaad674699bb607d8fedf1c2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27140b88353bd0c13c367205d9e1d43db25ace5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0b837e76f278922680c4f162175a90d15a1c438c/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fd7690f6242f2ab8aaad674699bb607d8fedf1c2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-28 09:51:25"
$ws3.Range("E2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$ws3.Range("F2").Value = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf"
$ws3.Range("G2").Value = "2016-01-28 09:52:25"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf"
$ws3.Range("D3").Value = "2016-01-28 09:51:25"
$ws3.Range("E3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$ws3.Range("F3").Value = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf"
$ws3.Range("G3").Value = "2016-01-28 09:52:25"
$ws3.Range("H3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3be81cff3d31c4421371fa4b241cfa5e6b3fbbe8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/66ae5c840647fb956586a03149f3c44423def458/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3483816a9b860720a66c4bdb7dba45f96721f277/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf", "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3be81cff3d31c4421371fa4b241cfa5e6b3fbbe8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/66ae5c840647fb956586a03149f3c44423def458/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3483816a9b860720a66c4bdb7dba45f96721f277/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf", "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c74a380b7562dbd8d9b487d73f16ef0a101212b9/.localization-config", "", "", ".localization-config") | Out-Null
